$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "Dinafex 180mg Tablet"
$ws.Range("C4").Value = "Dinafex 120mg Tablet"
$ws.Range("C5").Value = "Dinafex 60mg Tablet"

$ws.Range("C7").Value = "Etorix 120mg Tablet"
$ws.Range("D7").Value = "20's"
$ws.Range("C8").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("D8").Value = "40's"
$ws.Range("C9").Value = "Etorix 90mg Tablet"
$ws.Range("D9").Value = "30's"

$ws.Range("C11").Value = "Flucloxin 500mg Capsule"
$ws.Range("D11").Value = "30 's"
$ws.Range("C12").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("D12").Value = "36 's"

$ws.Range("F13").Value = 310
$ws.Range("G13").Value = 31

$ws.Range("C15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("D15").Value = "4's"
$ws.Range("C16").Value = "Ketonic 10mg Tablet"
$ws.Range("D16").Value = "20's"

$ws.Range("C17").Value = "Kynol D 25mg Tablet"
$ws.Range("D17").Value = "60 's"
$ws.Range("C18").Value = "Kynol TR 100mg Capsule"
$ws.Range("D18").Value = "50 's"
$ws.Range("C19").Value = "Kynol TR 200mg Capsule"
$ws.Range("D19").Value = "30 's"

$ws.Range("C24").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("D24").Value = "30ml"
$ws.Range("C25").Value = "Zithrox 15ml Suspension"
$ws.Range("D25").Value = "15 ml"
$ws.Range("C27").Value = "Zithrox 500mg Tablet"
$ws.Range("D27").Value = "6 's"

Write-Output "Done applying branch wise stock edits"
